# "Manutenção no artefato 15"
#
# 1. Remove the whole "Cenário: Fornecer serviço de guia" slide (was slide 5).
# 2. Remove the three shapes on slide 1 that referenced that scenario
#    (Retângulo 9 / Conector reto 11 / CaixaDeTexto 13).
# 3. Re-create the default PowerPoint section over the remaining 4 slides.
# 4. Refresh the cached "datetimeFigureOut" footer field from 08/03/2021 to
#    11/03/2021 on every slide layout and on the slide master.

$p = $ppt.ActivePresentation

# --- 1. Remove the shapes tied to the scenario that is going away ---------
$s1 = $p.Slides.Item(1)
for ($i = $s1.Shapes.Count; $i -ge 1; $i--) {
    $shape = $s1.Shapes.Item($i)
    if ($shape.Id -eq 10 -or $shape.Id -eq 12 -or $shape.Id -eq 14) {
        $shape.Delete()
    }
}

# --- 2. Delete the now-orphaned "Fornecer serviço de guia" slide ----------
$p.Slides.Item($p.Slides.Count).Delete()

# --- 3. Recreate the default section across the remaining slides ----------
$p.SectionProperties.AddSection(1, "Seção Padrão") | Out-Null

# ppPlaceholderDate = 16 : works regardless of the UI locale (shape Names
# are the Portuguese "Espaço Reservado para Data N" in this deck).
$ppPlaceholderDate = 16

function Update-DatePlaceholder($shapes) {
    foreach ($shape in $shapes) {
        if (-not $shape.HasTextFrame) { continue }
        $isDatePh = $false
        try {
            if ($shape.PlaceholderFormat.Type -eq $ppPlaceholderDate) { $isDatePh = $true }
        } catch {
            $isDatePh = $false
        }
        if ($isDatePh -and $shape.TextFrame.TextRange.Text -eq "08/03/2021") {
            $shape.TextFrame.TextRange.Text = "11/03/2021"
        }
    }
}

$master = $p.SlideMaster
Update-DatePlaceholder $master.Shapes

foreach ($layout in $master.CustomLayouts) {
    Update-DatePlaceholder $layout.Shapes
}

Write-Output "done"
